# Insert a new weekly price record as row 105 (pushing the existing
# rows 105-131 down to 106-132), for "Fruta / Arándano (blue)" at
# "Vega Modelo de Temuco".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 105; Excel shifts rows 105..131 down to 106..132
# and the sheet dimension grows from A1:T131 to A1:T132 automatically.
$ws.Rows.Item(105).Insert()

# Populate the newly inserted row 105 with the new weekly record.
$ws.Range("A105").Value = 10
$ws.Range("B105").Value = "Vega Modelo de Temuco"
$ws.Range("C105").Value = "La Araucanía"
$ws.Range("D105").Value = 44932
$ws.Range("E105").Value = 9
$ws.Range("F105").Value = "Fruta"
$ws.Range("G105").Value = 100101
$ws.Range("H105").Value = "Berries"
$ws.Range("I105").Value = 100101001
$ws.Range("J105").Value = "Arándano (blue)"
$ws.Range("K105").Value = "Sin especificar"
$ws.Range("L105").Value = "Primera"
$ws.Range("M105").Value = 125
$ws.Range("N105").Value = 1800
$ws.Range("O105").Value = 1800
$ws.Range("P105").Value = 1800
$ws.Range("Q105").Value = "$/kilo"
$ws.Range("R105").Value = "Región del Maule"
$ws.Range("S105").Value = 1800
$ws.Range("T105").Value = 1
